$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.770638372978055
$ws.Range("C2").Value = 0.01311926495731797
$ws.Range("D2").Value = 0.03223974428588861
$ws.Range("E2").Value = 0.06383820745116875
$ws.Range("F2").Value = 7.662150779775487
$ws.Range("H2").Value = 0.07973214163530429
$ws.Range("J2").Value = 0.2588048803500271
$ws.Range("K2").Value = 1.236897356269907
$ws.Range("L2").Value = 0.2246228339346352
$ws.Range("M2").Value = 0.3682373811667894

$ws.Range("B3").Value = 1.768299610546876
$ws.Range("C3").Value = 0.01138622289478519
$ws.Range("D3").Value = 0.02831459372016809
$ws.Range("E3").Value = 0.06421076739112808
$ws.Range("F3").Value = 7.468097427213934
$ws.Range("H3").Value = 0.07973214163530429
$ws.Range("J3").Value = 0.2557179444048998
$ws.Range("K3").Value = 1.231879661895604
$ws.Range("L3").Value = 0.226917329774956
$ws.Range("M3").Value = 0.3693921511879665

$ws.Range("B4").Value = 1.768133792704674
$ws.Range("C4").Value = 0.01033676375431725
$ws.Range("D4").Value = 0.02589476758429043
$ws.Range("E4").Value = 0.06445473036807403
$ws.Range("F4").Value = 7.349342984866212
$ws.Range("H4").Value = 0.07973214163530429
$ws.Range("J4").Value = 0.2538144995236209
$ws.Range("K4").Value = 1.229838432415775
$ws.Range("L4").Value = 0.2284582534633905
$ws.Range("M4").Value = 0.3703540916250638

$ws.Range("B5").Value = 1.768385646666701
$ws.Range("C5").Value = 0.0099126790523556
$ws.Range("D5").Value = 0.02490607137514189
$ws.Range("E5").Value = 0.06455798277131741
$ws.Range("F5").Value = 7.301045931886136
$ws.Range("H5").Value = 0.07973214163530429
$ws.Range("J5").Value = 0.2530367374713052
$ws.Range("K5").Value = 1.229267871343922
$ws.Range("L5").Value = 0.2291194510684456
$ws.Range("M5").Value = 0.3708097158614052

$ws.Range("B6").Value = 1.768446760421824
$ws.Range("C6").Value = 0.009842473273252494
$ws.Range("D6").Value = 0.02474173811978631
$ws.Range("E6").Value = 0.06457535973083672
$ws.Range("F6").Value = 7.293031974943716
$ws.Range("H6").Value = 0.07973214163530429
$ws.Range("J6").Value = 0.2529074623056289
$ws.Range("K6").Value = 1.229188904547087
$ws.Range("L6").Value = 0.2292312524360227
$ws.Range("M6").Value = 0.3708892156880026

$ws.Range("B7").Value = 1.768135895918107
$ws.Range("C7").Value = 0.01033103004533586
$ws.Range("D7").Value = 0.02588144438298912
$ws.Range("E7").Value = 0.06445610732184548
$ws.Range("F7").Value = 7.348691247812241
$ws.Range("H7").Value = 0.07973214163530429
$ws.Range("J7").Value = 0.2538040189134847
$ws.Range("K7").Value = 1.22982968004176
$ws.Range("L7").Value = 0.2284670358859771
$ws.Range("M7").Value = 0.3703599786753315

$ws.Range("B8").Value = 1.769568343775376
$ws.Range("C8").Value = 0.01251861243860475
$ws.Range("D8").Value = 0.03088829371068869
$ws.Range("E8").Value = 0.06396351591002958
$ws.Range("F8").Value = 7.595157307588408
$ws.Range("H8").Value = 0.07973214163530429
$ws.Range("J8").Value = 0.2577421348041753
$ws.Range("K8").Value = 1.234951360355609
$ws.Range("L8").Value = 0.2253865852762011
$ws.Range("M8").Value = 0.3685830714677429

$ws.Range("B9").Value = 1.782458746692697
$ws.Range("C9").Value = 0.01692939683639594
$ws.Range("D9").Value = 0.04063595127169606
$ws.Range("E9").Value = 0.06311773911261653
$ws.Range("F9").Value = 8.081771553619461
$ws.Range("H9").Value = 0.07973214163530429
$ws.Range("J9").Value = 0.2654042491397917
$ws.Range("K9").Value = 1.25325567135198
$ws.Range("L9").Value = 0.2203921392918602
$ws.Range("M9").Value = 0.3671045329616156

$ws.Range("B10").Value = 1.798085353570116
$ws.Range("C10").Value = 0.02025071311219051
$ws.Range("D10").Value = 0.04776429543568383
$ws.Range("E10").Value = 0.0625689500965958
$ws.Range("F10").Value = 8.441552870573929
$ws.Range("H10").Value = 0.07973214163530429
$ws.Range("J10").Value = 0.2710018749576903
$ws.Range("K10").Value = 1.27176111505571
$ws.Range("L10").Value = 0.2173582705301769
$ws.Range("M10").Value = 0.367241047455952

$ws.Range("B11").Value = 1.806533849416496
$ws.Range("C11").Value = 0.02178072871079451
$ws.Range("D11").Value = 0.05100210241059244
$ws.Range("E11").Value = 0.06233491532618096
$ws.Range("F11").Value = 8.60577662053413
$ws.Range("H11").Value = 0.07973214163530429
$ws.Range("J11").Value = 0.2735426346901448
$ws.Range("K11").Value = 1.281283015418893
$ws.Range("L11").Value = 0.21611562633489
$ws.Range("M11").Value = 0.3675687559826208

$ws.Range("B12").Value = 1.809925905342936
$ws.Range("C12").Value = 0.02236297684788724
$ws.Range("D12").Value = 0.05222763200272595
$ws.Range("E12").Value = 0.06224852647133305
$ws.Range("F12").Value = 8.668048006292338
$ws.Range("H12").Value = 0.07973214163530429
$ws.Range("J12").Value = 0.2745040272888701
$ws.Range("K12").Value = 1.285047772409342
$ws.Range("L12").Value = 0.2156648027488615
$ws.Range("M12").Value = 0.3677310434028165

$ws.Range("B13").Value = 1.809186788744796
$ws.Range("C13").Value = 0.02223745008178923
$ws.Range("D13").Value = 0.05196371470429995
$ws.Range("E13").Value = 0.062267032611107
$ws.Range("F13").Value = 8.654632986926401
$ws.Range("H13").Value = 0.07973214163530429
$ws.Range("J13").Value = 0.2742970058711549
$ws.Range("K13").Value = 1.284229887053016
$ws.Range("L13").Value = 0.215761018257389
$ws.Range("M13").Value = 0.3676943933480743

$ws.Range("B14").Value = 1.806809051271358
$ws.Range("C14").Value = 0.02182857255530735
$ws.Range("D14").Value = 0.05110293778919583
$ws.Range("E14").Value = 0.06232776332653156
$ws.Range("F14").Value = 8.610898034335889
$ws.Range("H14").Value = 0.07973214163530429
$ws.Range("J14").Value = 0.273621743416939
$ws.Range("K14").Value = 1.281589555407407
$ws.Range("L14").Value = 0.2160781413824537
$ws.Range("M14").Value = 0.3675813420615732

$ws.Range("B15").Value = 1.805377729740485
$ws.Range("C15").Value = 0.02157849968112657
$ws.Range("D15").Value = 0.05057561887986139
$ws.Range("E15").Value = 0.06236525341612698
$ws.Range("F15").Value = 8.584120073638189
$ws.Range("H15").Value = 0.07973214163530429
$ws.Range("J15").Value = 0.2732080316886965
$ws.Range("K15").Value = 1.27999299463346
$ws.Range("L15").Value = 0.2162749581448935
$ws.Range("M15").Value = 0.3675170684862543

$ws.Range("B16").Value = 1.797560196240511
$ws.Range("C16").Value = 0.02015111763624589
$ws.Range("D16").Value = 0.04755261098659957
$ws.Range("E16").Value = 0.06258455812071517
$ws.Range("F16").Value = 8.430831956774682
$ws.Range("H16").Value = 0.07973214163530429
$ws.Range("J16").Value = 0.2708357231669751
$ws.Range("K16").Value = 1.271161071702721
$ws.Range("L16").Value = 0.21744224403907
$ws.Range("M16").Value = 0.3672249753398482

$ws.Range("B17").Value = 1.793107648253596
$ws.Range("C17").Value = 0.01928044539631912
$ws.Range("D17").Value = 0.04569694100287336
$ws.Range("E17").Value = 0.06272308596978249
$ws.Range("F17").Value = 8.336939738165938
$ws.Range("H17").Value = 0.07973214163530429
$ws.Range("J17").Value = 0.2693789978016099
$ws.Range("K17").Value = 1.266025869729674
$ws.Range("L17").Value = 0.2181935255838141
$ws.Range("M17").Value = 0.3671138112972372

$ws.Range("B18").Value = 1.790672767517151
$ws.Range("C18").Value = 0.01878145171926349
$ws.Range("D18").Value = 0.04462913738453267
$ws.Range("E18").Value = 0.06280423360122445
$ws.Range("F18").Value = 8.282987687559626
$ws.Range("H18").Value = 0.07973214163530429
$ws.Range("J18").Value = 0.2685405938360432
$ws.Range("K18").Value = 1.263176101569258
$ws.Range("L18").Value = 0.2186385846830632
$ws.Range("M18").Value = 0.367074872906592

$ws.Range("B19").Value = 1.789870014279984
$ws.Range("C19").Value = 0.01861280571770862
$ws.Range("D19").Value = 0.04426751195460099
$ws.Range("E19").Value = 0.06283196158679782
$ws.Range("F19").Value = 8.264729329039255
$ws.Range("H19").Value = 0.07973214163530429
$ws.Range("J19").Value = 0.2682566307676524
$ws.Range("K19").Value = 1.262229048527075
$ws.Range("L19").Value = 0.2187914977587226
$ws.Range("M19").Value = 0.3670659829732372

$ws.Range("B20").Value = 1.793568577180025
$ws.Range("C20").Value = 0.01937294345638918
$ws.Range("D20").Value = 0.04589452794172644
$ws.Range("E20").Value = 0.06270818736912709
$ws.Range("F20").Value = 8.346929291915018
$ws.Range("H20").Value = 0.07973214163530429
$ws.Range("J20").Value = 0.2695341234583353
$ws.Range("K20").Value = 1.266561769006671
$ws.Range("L20").Value = 0.2181122111773348
$ws.Range("M20").Value = 0.3671230575052995

$ws.Range("B21").Value = 1.807502216987473
$ws.Range("C21").Value = 0.0219485911478472
$ws.Range("D21").Value = 0.05135578278006392
$ws.Range("E21").Value = 0.06230986466413402
$ws.Range("F21").Value = 8.623741760012706
$ws.Range("H21").Value = 0.07973214163530429
$ws.Range("J21").Value = 0.2738201038253791
$ws.Range("K21").Value = 1.282360766500204
$ws.Range("L21").Value = 0.2159844591645523
$ws.Range("M21").Value = 0.3676135114715748

$ws.Range("B22").Value = 1.817732473885201
$ws.Range("C22").Value = 0.02364867535365534
$ws.Range("D22").Value = 0.05492184768651498
$ws.Range("E22").Value = 0.06206256066765459
$ws.Range("F22").Value = 8.805142461588048
$ws.Range("H22").Value = 0.07973214163530429
$ws.Range("J22").Value = 0.2766169718092186
$ws.Range("K22").Value = 1.29361336980881
$ws.Range("L22").Value = 0.2147088897312699
$ws.Range("M22").Value = 0.3681566603737991

$ws.Range("B23").Value = 1.812169519645181
$ws.Range("C23").Value = 0.02273973809657548
$ws.Range("D23").Value = 0.05301881345295101
$ws.Range("E23").Value = 0.06219336313932811
$ws.Range("F23").Value = 8.708279776794654
$ws.Range("H23").Value = 0.07973214163530429
$ws.Range("J23").Value = 0.2751245975259664
$ws.Range("K23").Value = 1.287522712131334
$ws.Range("L23").Value = 0.2153791687733317
$ws.Range("M23").Value = 0.3678464030391595

$ws.Range("B24").Value = 1.793359802020149
$ws.Range("C24").Value = 0.01933112021048089
$ws.Range("D24").Value = 0.04580520187191439
$ws.Range("E24").Value = 0.06271491833438692
$ws.Range("F24").Value = 8.342412928184785
$ws.Range("H24").Value = 0.07973214163530429
$ws.Range("J24").Value = 0.2694639940007661
$ws.Range("K24").Value = 1.266319169603662
$ws.Range("L24").Value = 0.2181489324969661
$ws.Range("M24").Value = 0.3671187995112746

$ws.Range("B25").Value = 1.777891274930028
$ws.Range("C25").Value = 0.01572249728007336
$ws.Range("D25").Value = 0.03800550784993106
$ws.Range("E25").Value = 0.06333374616470433
$ws.Range("F25").Value = 7.949750310324447
$ws.Range("H25").Value = 0.07973214163530429
$ws.Range("J25").Value = 0.2633373910601122
$ws.Range("K25").Value = 1.24741735091763
$ws.Range("L25").Value = 0.2216314970001747
$ws.Range("M25").Value = 0.3672897833950657
